$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted at the top of the
# "Femacal de La Calera - Arándano (blue)" data block (row 317),
# pushing all the existing rows 317-380 down to 318-381.
$ws.Rows("317:317").Insert()

$ws.Range("A317").Value = 3
$ws.Range("B317").Value = "Femacal de La Calera"
$ws.Range("C317").Value = "Coquimbo"
$ws.Range("D317").Value = 45211
$ws.Range("E317").Value = 5
$ws.Range("F317").Value = "Fruta"
$ws.Range("G317").Value = 100101
$ws.Range("H317").Value = "Berries"
$ws.Range("I317").Value = 100101001
$ws.Range("J317").Value = "Arándano (blue)"
$ws.Range("K317").Value = "Sin especificar"
$ws.Range("L317").Value = "Primera"
$ws.Range("M317").Value = 75
$ws.Range("N317").Value = 12000
$ws.Range("O317").Value = 13000
$ws.Range("P317").Value = 12600
$ws.Range("Q317").Value = "$/bandeja 2 kilos"
$ws.Range("R317").Value = "Provincia de Quillota"
$ws.Range("S317").Value = 6300
$ws.Range("T317").Value = 2
